$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 16

$ws.Cells.Item($row, 1).Value = 7
$ws.Cells.Item($row, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item($row, 3).Value = 'Ñuble'
$ws.Cells.Item($row, 4).Value = 45194
$ws.Cells.Item($row, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item($row, 5).Value = 16
$ws.Cells.Item($row, 6).Value = 'Fruta'
$ws.Cells.Item($row, 7).Value = 100102
$ws.Cells.Item($row, 8).Value = 'Cítricos'
$ws.Cells.Item($row, 9).Value = 100102006
$ws.Cells.Item($row, 10).Value = 'Pomelo'
$ws.Cells.Item($row, 11).Value = 'Start Ruby'
$ws.Cells.Item($row, 12).Value = 'Primera'
$ws.Cells.Item($row, 13).Value = 60
$ws.Cells.Item($row, 14).Value = 15000
$ws.Cells.Item($row, 15).Value = 15000
$ws.Cells.Item($row, 16).Value = 15000
$ws.Cells.Item($row, 17).Value = '$/caja 14 kilos granel'
$ws.Cells.Item($row, 18).Value = "Región de O'Higgins"
$ws.Cells.Item($row, 19).Value = 1071
$ws.Cells.Item($row, 20).Value = 14
